$wb = $excel.ActiveWorkbook

# Sheet ALC, row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3004
$ws.Range("J38").Value = 4000
$ws.Range("L38").Value = 12000
$ws.Range("N38").Value = -12744

# Sheet ALC, row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 5694.4614
$ws.Range("I58").Value = 286
$ws.Range("K58").Value = 858
$ws.Range("M58").Value = -708

# Sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3529.2
$ws.Range("I70").Value = 2548.6667
$ws.Range("K70").Value = 7646.000100000001
$ws.Range("M70").Value = -7376.000100000001

# Sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3529.2
$ws.Range("I73").Value = 2548.6667
$ws.Range("K73").Value = 7646.000100000001
$ws.Range("M73").Value = -6710.000100000001

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 22729416
$ws.Range("J107").Value = 2397.4443
$ws.Range("L107").Value = 2397.4443
$ws.Range("N107").Value = -6237.4443

# Sheet ALC, row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 51306.668
$ws.Range("J133").Value = 51306.668
$ws.Range("L133").Value = 51306.668
$ws.Range("N133").Value = -61426.668

# Sheet ALC, row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 56808.89
$ws.Range("J136").Value = 56808.89
$ws.Range("L136").Value = 56808.89
$ws.Range("N136").Value = -67008.89

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4217.2393
$ws.Range("I138").Value = 2355.6667
$ws.Range("J138").Value = 4595.8643
$ws.Range("K138").Value = 7067.000100000001
$ws.Range("L138").Value = 13787.5929
$ws.Range("M138").Value = -1927.000100000001
$ws.Range("N138").Value = -24067.5929

# Sheet ALC, row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 56623.332
$ws.Range("J139").Value = 56623.332
$ws.Range("L139").Value = 56623.332
$ws.Range("N139").Value = -66903.33199999999

# Sheet ARM, row 106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 43000
$ws.Range("J106").Value = 43000
$ws.Range("L106").Value = 43000
$ws.Range("N106").Value = -45524

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 9289.137000000001
$ws.Range("I110").Value = 10932.206
$ws.Range("J110").Value = 3702.7
$ws.Range("K110").Value = 10932.206
$ws.Range("L110").Value = 3702.7
$ws.Range("M110").Value = -8887.206
$ws.Range("N110").Value = -7792.7

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2980.149
$ws.Range("I132").Value = 1871.5161
$ws.Range("J132").Value = 5128.125
$ws.Range("K132").Value = 5614.5483
$ws.Range("L132").Value = 15384.375
$ws.Range("M132").Value = -3084.5483
$ws.Range("N132").Value = -20444.375

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2200
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2200
$ws.Range("N86").Value = -4446
$ws.Range("M86").Value = $null

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2200
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 11000
$ws.Range("N89").Value = -22232
$ws.Range("M89").Value = $null

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3000011
$ws.Range("I107").Value = 3000011
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3000011
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = -2998091
$ws.Range("M107").Value = $null

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15155410
$ws.Range("I31").Value = 1699.48
$ws.Range("J31").Value = 62510750
$ws.Range("K31").Value = 1699.48
$ws.Range("L31").Value = 62510750
$ws.Range("M31").Value = -1404.48
$ws.Range("N31").Value = -62511340

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 15155410
$ws.Range("I34").Value = 1699.48
$ws.Range("J34").Value = 62510750
$ws.Range("K34").Value = 1699.48
$ws.Range("L34").Value = 62510750
$ws.Range("M34").Value = -1497.48
$ws.Range("N34").Value = -62511154

# Sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6110
$ws.Range("I94").Value = 5313
$ws.Range("J94").Value = 6337.7144
$ws.Range("K94").Value = 5313
$ws.Range("L94").Value = 6337.7144
$ws.Range("M94").Value = -4862
$ws.Range("N94").Value = -7239.7144

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1299.4736
$ws.Range("I68").Value = 674.375
$ws.Range("J68").Value = 4633.3335
$ws.Range("K68").Value = 2023.125
$ws.Range("L68").Value = 13900.0005
$ws.Range("M68").Value = -1212.125
$ws.Range("N68").Value = -15522.0005

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1299.4736
$ws.Range("I71").Value = 674.375
$ws.Range("J71").Value = 4633.3335
$ws.Range("K71").Value = 6069.375
$ws.Range("L71").Value = 41700.0015
$ws.Range("M71").Value = -2013.375
$ws.Range("N71").Value = -49812.0015

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6928.3335
$ws.Range("I122").Value = 1127.8
$ws.Range("J122").Value = 14179
$ws.Range("K122").Value = 10150.2
$ws.Range("L122").Value = 127611
$ws.Range("M122").Value = -7700.199999999999
$ws.Range("N122").Value = -132511

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2326812.8
$ws.Range("I131").Value = 9091697
$ws.Range("J131").Value = 1383.6875
$ws.Range("K131").Value = 27275091
$ws.Range("L131").Value = 4151.0625
$ws.Range("M131").Value = -27270051
$ws.Range("N131").Value = -14231.0625

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 59341456
$ws.Range("I122").Value = 62636980
$ws.Range("J122").Value = 50004136
$ws.Range("K122").Value = 187910940
$ws.Range("L122").Value = 150012408
$ws.Range("M122").Value = -187908490
$ws.Range("N122").Value = -150017308

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6946704
$ws.Range("I132").Value = 8773679
$ws.Range("J132").Value = 4199.8
$ws.Range("K132").Value = 26321037
$ws.Range("L132").Value = 12599.4
$ws.Range("M132").Value = -26318507
$ws.Range("N132").Value = -17659.4

# Sheet LTW, row 35
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1663.1428
$ws.Range("I35").Value = 940.3333
$ws.Range("K35").Value = 940.3333
$ws.Range("M35").Value = -604.3333

# Sheet LTW, row 121
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null

# Sheet WVR, row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 7213.5
$ws.Range("J21").Value = 7213.5
$ws.Range("L21").Value = 7213.5
$ws.Range("N21").Value = -7683.5

# Sheet WVR, row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 7213.5
$ws.Range("J35").Value = 7213.5
$ws.Range("L35").Value = 7213.5
$ws.Range("N35").Value = -7793.5

# Sheet WVR, row 93
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2297.3809
$ws.Range("I132").Value = 976.2
$ws.Range("J132").Value = 3498.4546
$ws.Range("K132").Value = 2928.6
$ws.Range("L132").Value = 10495.3638
$ws.Range("M132").Value = -398.6000000000004
$ws.Range("N132").Value = -15555.3638
